# Updates cryptos list price/volume figures (scheduled data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.621.76'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '1.860.16'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.21'
$ws.Range('E5').Value = '  +2.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6976'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07717'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3063'
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.63'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07755'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.161'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '1.856.24'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.25'
$ws.Range('E14').Value = '  +2.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6922'
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.554'
$ws.Range('E16').Value = '  +3.09%  '
$ws.Range('D17').Value = '29.610.37'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008342'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').Value = '2.105.38'
$ws.Range('E19').Value = '  +1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.67'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.75'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.608'
$ws.Range('E23').Value = '  +2.87%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1505'
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.909'
$ws.Range('E26').Value = '  +2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.67'
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.29'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.538'
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.250'
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.193'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05099'
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7787'
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.898'
$ws.Range('E35').Value = '  +5.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.154'
$ws.Range('E36').Value = '  +1.60%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').Value = '1.326.23'
$ws.Range('E38').Value = '  +10.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01874'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.734'
$ws.Range('E40').Value = '  +2.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9588'
$ws.Range('E41').Value = '  +3.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.26'
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.840'
$ws.Range('E43').Value = '  +12.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9998'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('E45').Value = '  +4.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.769'
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('D47').Value = '2.005.09'
$ws.Range('E47').Value = '  +1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5213'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.783'
$ws.Range('E49').Value = '  +3.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '64.53'
$ws.Range('E50').Value = '  +4.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.984'
$ws.Range('E51').Value = '  +1.81%  '
